$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 data
$ws.Range("A3").Value = "Tuesday"
$ws.Range("B3").Value = 45769
$ws.Range("C3").Value = 0.70833333333333337
$ws.Range("D3").Value = 0.79166666666666663
$ws.Range("E3").Value = 0.083333333333333329
$ws.Range("F3").Value = "Started working on character controls"

# Reuse the same number formats/styles as used on the row above,
# by copying formats only (so existing style indices are reused).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C3:D3").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
